# Journal de travail - ajout d'une entree de suivi pour la phase de
# placement/arrangement d'un regiment (analyse + debut d'explication de
# l'algorithme), et mise a jour de la vue/de la mise en page associees.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Donnees de la ligne 51 (Fin / Heure / Description) -------------------
# Fin (C51) : l'heure de fin de la tache a ete renseignee.
$ws.Range("C51").Value = 0.70486111111111116

# Description (F51) : note sur la difficulte a expliquer l'algorithme.
$ws.Range("F51").Value = "Expliquer un algorithme est plus long qu'imaginé, des sacrifice vont devoir être fait"

# --- Vue de la feuille : position de defilement et selection --------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 22
$ws.Range("A1:H51").Select()

# --- Mise en page : echelle d'impression -----------------------------------
$ps = $ws.PageSetup
$ps.Zoom = 33
$ps.FitToPagesTall = $false
